$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 810.6731951491925
$ws.Range("C2").Value = 519.5399500817239
$ws.Range("D2").Value = 377.5165848720561
$ws.Range("E2").Value = 891.2072356207559
$ws.Range("F2").Value = 816.2655584213176
$ws.Range("G2").Value = 705.698026459629
$ws.Range("H2").Value = 894.2834347060425

$ws.Range("B3").Value = 853.6162607233118
$ws.Range("C3").Value = 542.3516677527452
$ws.Range("D3").Value = 354.2488000318697
$ws.Range("E3").Value = 918.8464568727627
$ws.Range("F3").Value = 857.5566018388613
$ws.Range("G3").Value = 713.4871374663375
$ws.Range("H3").Value = 920.7716035501105

$ws.Range("B4").Value = 734.5948601966046
$ws.Range("C4").Value = 563.8605507489592
$ws.Range("D4").Value = 200.3093732408394
$ws.Range("E4").Value = 848.1439060241668
$ws.Range("F4").Value = 741.6630058926968
$ws.Range("G4").Value = 636.8114506417281
$ws.Range("H4").Value = 851.6346524056777

$ws.Range("B5").Value = 915.9270000733128
$ws.Range("C5").Value = 662.9786698217429
$ws.Range("D5").Value = 130.3464593441165
$ws.Range("E5").Value = 965.1838895242677
$ws.Range("F5").Value = 917.083831534653
$ws.Range("G5").Value = 714.0936209589279
$ws.Range("H5").Value = 965.6008698726566

$ws.Range("B6").Value = 900.2329883775661
$ws.Range("C6").Value = 647.8814483192817
$ws.Range("D6").Value = 66.49382152439483
$ws.Range("E6").Value = 956.0626579123599
$ws.Range("F6").Value = 901.9208893662466
$ws.Range("G6").Value = 671.7520134723012
$ws.Range("H6").Value = 956.6387882867735

$ws.Range("B7").Value = 938.8517375560593
$ws.Range("C7").Value = 713.9659208237927
$ws.Range("D7").Value = 111.4180642391784
$ws.Range("E7").Value = 978.3911234606398
$ws.Range("F7").Value = 939.6353147012608
$ws.Range("G7").Value = 750.8099286789158
$ws.Range("H7").Value = 978.6812250594128

$ws.Range("B8").Value = 833.1799429183484
$ws.Range("C8").Value = 285.6775235257514
$ws.Range("D8").Value = 108.5610867388323
$ws.Range("E8").Value = 857.551967363283
$ws.Range("F8").Value = 833.2110855707815
$ws.Range("G8").Value = 356.6613353955055
$ws.Range("H8").Value = 857.5831100157161

$ws.Range("B9").Value = 945.7389134525196
$ws.Range("C9").Value = 606.6176074155284
$ws.Range("D9").Value = 175.6678320388319
$ws.Range("E9").Value = 968.2246803912416
$ws.Range("F9").Value = 946.7832509100929
$ws.Range("G9").Value = 665.5488374974055
$ws.Range("H9").Value = 968.8104278695378

$ws.Range("B10").Value = 961.3154436913049
$ws.Range("C10").Value = 683.4259752372587
$ws.Range("D10").Value = 266.0637605274405
$ws.Range("E10").Value = 984.3029141654705
$ws.Range("F10").Value = 962.4164304111293
$ws.Range("G10").Value = 796.7028593535204
$ws.Range("H10").Value = 985.2174112199041

$ws.Range("B11").Value = 952.5553249171585
$ws.Range("C11").Value = 686.6652809504774
$ws.Range("D11").Value = 232.4492655713701
$ws.Range("E11").Value = 983.7696678525957
$ws.Range("F11").Value = 953.6963681437817
$ws.Range("G11").Value = 788.236008520571
$ws.Range("H11").Value = 984.5733410625709

$ws.Range("B12").Value = 874.7493301769194
$ws.Range("C12").Value = 700.8260013260174
$ws.Range("D12").Value = 56.26714088105341
$ws.Range("E12").Value = 962.237953507247
$ws.Range("F12").Value = 877.6914250533716
$ws.Range("G12").Value = 720.6201348792119
$ws.Range("H12").Value = 963.1117077385796

$ws.Range("B13").Value = 956.7539180159872
$ws.Range("C13").Value = 681.6785216095993
$ws.Range("D13").Value = 203.0457098268768
$ws.Range("E13").Value = 982.220056285759
$ws.Range("F13").Value = 957.8519584561764
$ws.Range("G13").Value = 774.2730251325225
$ws.Range("H13").Value = 983.0675815936149

